$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new columns at G:K (pushes existing G:V -> L:AA),
# carrying along existing per-cell formatting via the native Insert shift.
$ws.Range("G1:K1").EntireColumn.Insert()

# --- Header row (row 1): new "meta" group headers ---
$ws.Range("G1").Value = "meta"
$ws.Range("H1").Value = "meta_avg"
$ws.Range("I1").Value = "meta_std"
$ws.Range("J1").Value = "meta_min"
$ws.Range("K1").Value = "meta_max"

# --- Data row (row 2): new "meta" values, currency-formatted like the
#     neighbouring arrecadado_* columns ---
$ws.Range("G2:K2").NumberFormat = "R$ #,##0.00"
$ws.Range("G2").Value = 13973042.60019265
$ws.Range("H2").Value = 16834.99108456945
$ws.Range("I2").Value = 17015.69760983049
$ws.Range("J2").Value = 31.89582864100442
$ws.Range("K2").Value = 189313.7035611726

# --- Updated statistic: contribuicoes_std (now column W after the shift) ---
$ws.Range("W2").Value = 423.0192251466749
